$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Q0-Q7 (rows 37-44): rotation 270 -> 90
for ($r = 37; $r -le 44; $r++) {
    $ws.Range("E$r").Value = 90
}

# Q8-Q15 (rows 45-52): rotation 90 -> 270, and add an (empty) formatted cell in column F
for ($r = 45; $r -le 52; $r++) {
    $ws.Range("E$r").Value = 270
    $ws.Range("F$r").Style = $ws.Range("E$r").Style
}

# Rows 143-149: update rotation values
$ws.Range("E143").Value = 270
$ws.Range("E144").Value = 90
$ws.Range("E145").Value = 270
$ws.Range("E146").Value = 270
$ws.Range("E147").Value = 270
$ws.Range("E148").Value = 270
$ws.Range("E149").Value = 270

# Update the active selection to match the committed view state
$ws.Range("F143").Select()
